$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '58.140.86'
$ws.Range('E2').Value = '  -3.49%  '
$ws.Range('D3').Value = '3.130.27'
$ws.Range('E3').Value = '  -5.29%  '
$ws.Range('E4').Value = '  +0.02%  '
$ws.Range('D5').Value = "'523.91"
$ws.Range('E5').Value = '  -6.04%  '
$ws.Range('D6').Value = "'134.23"
$ws.Range('E6').Value = '  -5.21%  '
$ws.Range('E7').Value = '  -0.06%  '
$ws.Range('D8').Value = '3.130.53'
$ws.Range('E8').Value = '  -5.27%  '
$ws.Range('D9').Value = "'0.439"
$ws.Range('E9').Value = '  -5.93%  '
$ws.Range('E10').Value = '  -8.29%  '
$ws.Range('E11').Value = '  -9.07%  '
$ws.Range('D12').Value = "'0.377"
$ws.Range('E12').Value = '  -7.70%  '
$ws.Range('D13').Value = '3.666.57'
$ws.Range('E13').Value = '  -5.26%  '
$ws.Range('E14').Value = '  -1.71%  '
$ws.Range('D15').Value = "'25.44"
$ws.Range('E15').Value = '  -5.56%  '
$ws.Range('D16').Value = '3.131.79'
$ws.Range('E16').Value = '  -5.12%  '
$ws.Range('D17').Value = '58.113.19'
$ws.Range('E17').Value = '  -3.63%  '
$ws.Range('D18').Value = "'0.0000152"
$ws.Range('E18').Value = '  -8.01%  '
$ws.Range('D19').Value = "'5.77"
$ws.Range('E19').Value = '  -5.52%  '
$ws.Range('D20').Value = "'13.01"
$ws.Range('E20').Value = '  -6.84%  '
$ws.Range('D21').Value = "'7.89"
$ws.Range('E21').Value = '  -7.92%  '
$ws.Range('D22').Value = "'342.43"
$ws.Range('E22').Value = '  -8.45%  '
$ws.Range('D23').Value = "'1.00"
$ws.Range('E23').Value = '  -0.14%  '
$ws.Range('B24').Value = 'Polygon'
$ws.Range('C24').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D24').Value = "'0.504"
$ws.Range('E24').Value = '  -5.52%  '
$ws.Range('B25').Value = 'Litecoin'
$ws.Range('C25').Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range('D25').Value = "'67.56"
$ws.Range('E25').Value = '  -8.95%  '
$ws.Range('D26').Value = '3.261.25'
$ws.Range('E26').Value = '  -5.12%  '
$ws.Range('B27').Value = 'PEPE'
$ws.Range('C27').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D27').Value = '0.0₃0955'
$ws.Range('E27').Value = '  -6.11%  '
$ws.Range('B28').Value = 'Kaspa'
$ws.Range('C28').Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range('D28').Value = "'0.167"
$ws.Range('E28').Value = '  -2.50%  '
$ws.Range('D29').Value = "'0.988"
$ws.Range('E29').Value = '  -1.10%  '
$ws.Range('E30').Value = '  -5.40%  '
$ws.Range('E32').Value = '  -8.63%  '
$ws.Range('D33').Value = "'6.81"
$ws.Range('E33').Value = '  -9.42%  '
$ws.Range('D34').Value = "'21.38"
$ws.Range('E34').Value = '  -5.46%  '
$ws.Range('D35').Value = "'1.21"
$ws.Range('E35').Value = '  -2.08%  '
$ws.Range('D36').Value = "'4.77"
$ws.Range('E36').Value = '  -6.34%  '
$ws.Range('D37').Value = "'156.54"
$ws.Range('E37').Value = '  -5.87%  '
$ws.Range('D38').Value = "'6.21"
$ws.Range('E38').Value = '  -6.80%  '
$ws.Range('E39').Value = '  -10.57%  '
$ws.Range('D40').Value = "'0.0684"
$ws.Range('E40').Value = '  -6.00%  '
$ws.Range('D41').Value = '3.160.89'
$ws.Range('E41').Value = '  -5.24%  '
$ws.Range('E42').Value = '  -3.89%  '
$ws.Range('D43').Value = "'24.10"
$ws.Range('E43').Value = '  -8.47%  '
$ws.Range('D44').Value = "'0.690"
$ws.Range('E44').Value = '  -8.03%  '
$ws.Range('E45').Value = '  -2.96%  '
$ws.Range('D46').Value = "'3.88"
$ws.Range('E46').Value = '  -5.87%  '
$ws.Range('E47').Value = '  -0.04%  '
$ws.Range('B48').Value = 'Stacks'
$ws.Range('C48').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D48').Value = "'1.43"
$ws.Range('E48').Value = '  -8.87%  '
$ws.Range('B49').Value = 'Maker'
$ws.Range('C49').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D49').Value = '2.258.00'
$ws.Range('E49').Value = '  -3.77%  '
$ws.Range('D50').Value = "'6.18"
$ws.Range('E50').Value = '  -3.02%  '
$ws.Range('D51').Value = "'20.59"
$ws.Range('E51').Value = '  -3.19%  '
